$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list (prices + 1h volume %) per the scheduled GitHub Actions refresh.

$ws.Range("D2").Value = '29.480.19'
$ws.Range("E2").Value = '  -0.03%  '

$ws.Range("D3").Value = '1.903.46'
$ws.Range("E3").Value = '  -0.26%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.45%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.55'
$ws.Range("E5").Value = '  -0.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.30%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4843'
$ws.Range("E7").Value = '  +3.96%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4061'
$ws.Range("E8").Value = '  -0.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08122'
$ws.Range("E9").Value = '  +1.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.007'
$ws.Range("E10").Value = '  +0.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.49'
$ws.Range("E11").Value = '  +5.45%  '

$ws.Range("D12").Value = '1.903.69'
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.016'
$ws.Range("E13").Value = '  +1.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.074'
$ws.Range("E14").Value = '  -0.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.37'
$ws.Range("E15").Value = '  +1.47%  '

$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06760'
$ws.Range("E17").Value = '  +2.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001039'
$ws.Range("E18").Value = '  +1.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.66'
$ws.Range("E19").Value = '  -0.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.38%  '

$ws.Range("D21").Value = '29.490.23'
$ws.Range("E21").Value = '  -0.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.564'
$ws.Range("E22").Value = '  +0.54%  '

$ws.Range("E23").Value = '  +2.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.160'
$ws.Range("E24").Value = '  -2.14%  '

$ws.Range("D25").Value = '2.140.97'
$ws.Range("E25").Value = '  +0.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.65'
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.06'
$ws.Range("E27").Value = '  +1.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.201'
$ws.Range("E28").Value = '  +8.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.098'
$ws.Range("E29").Value = '  -1.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.97'
$ws.Range("E30").Value = '  +1.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.031'
$ws.Range("E31").Value = '  -4.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09541'
$ws.Range("E32").Value = '  +0.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.518'
$ws.Range("E33").Value = '  +2.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.553'
$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.388'
$ws.Range("E35").Value = '  -2.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02261'
$ws.Range("E36").Value = '  +0.19%  '

$ws.Range("E37").Value = '  +0.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.172'
$ws.Range("E38").Value = '  +0.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5919'
$ws.Range("E39").Value = '  +0.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.905'
$ws.Range("E40").Value = '  -5.66%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.32'
$ws.Range("E41").Value = '  +2.37%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1853'
$ws.Range("E42").Value = '  +1.17%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.421'
$ws.Range("E43").Value = '  +1.97%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.284'
$ws.Range("E44").Value = '  -1.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07709'
$ws.Range("E45").Value = '  -0.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.34'
$ws.Range("E46").Value = '  +1.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5547'
$ws.Range("E47").Value = '  +0.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.939'
$ws.Range("E48").Value = '  +0.97%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.98'
$ws.Range("E49").Value = '  +1.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.60'
$ws.Range("E50").Value = '  +1.96%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.051'
$ws.Range("E51").Value = '  +2.30%  '
